# docs/protocol.xlsx: update to match reality
#
# Applies, via Excel COM automation, the same edits that were made by hand
# in the target commit:
#   1. Rewrite the long "helloData" description in F5 with the corrected
#      field-name/abbreviation mapping text (this also makes that row's
#      wrapped text taller).
#   2. Remove the two rows describing "tk_intraframe_corruption" and
#      "tk_brb, seconds" (rows 22 and 24) -- their row slots stay in place
#      (so every other row keeps its original row number) but all their
#      cell content/formatting is cleared.
#   3. Leave the selection on A5, scrolled so row 4 is at the top of the
#      window (mirrors the saved sheetView/selection state in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")
$ws.Activate()

# --- 1. Replace the helloData comment in F5 -------------------------------
$newHelloData = @"
 transportNumber = 'tnum', protocolVersion = 'ver', httpFormat = 'format', requestNewStream = 'new', streamId = 'id', credentialsData = 'cred',  streamingResponse = 'ming', needPaddingBytes = 'pad', maxReceiveBytes = 'maxb', maxOpenTime = 'maxt', useMyTcpAcks = 'tcpack', succeedsTransport = 'eeds', lastSackSeenByClient = 'lastack'
Presence of succeedsTransport option means "give me boxes, server". If succeedsTransport != null, temporarily assume that all boxes written to #<succeedsTransport> were SACKed.
Only C2S because hello is used by the client to identify itself to the server, and set critical transport parameters. XXX TODO: perhaps ackMode: 0 - require Minerva-level SACKs, 1 - use my TCP acks, 2 - assume everything written is received
"@

$ws.Range("F5").Value = $newHelloData

# Row 5 wraps text in F5, so the updated (longer) text needs a taller row.
$ws.Rows.Item(5).RowHeight = 115.5

# --- 2. Drop the "tk_intraframe_corruption" and "tk_brb, seconds" rows ----
$ws.Range("A22:F22").Clear()
$ws.Range("A24:F24").Clear()

# --- 3. Update the saved selection/scroll position -------------------------
$ws.Range("A5").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
